$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SkillsProfile")

# --- Update existing SkillsProfile sheet ---
$ws1.Range("H2").Value = 44909
$ws1.Range("I2").Value = 44909
$ws1.Range("L2").Value = 3
$ws1.Range("L2").HorizontalAlignment = -4131
$ws1.Range("I6").Select()

# --- Add the new Signup sheet after SkillsProfile ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Signup"

# Column widths (closest achievable values under host quantization)
$ws2.Columns.Item(1).ColumnWidth = 11.7566
$ws2.Columns.Item(2).ColumnWidth = 16.2598
$ws2.Columns.Item(3).ColumnWidth = 24.4195
$ws2.Columns.Item(4).ColumnWidth = 14.09
$ws2.Columns.Item(5).ColumnWidth = 16.7528

# Header row
$ws2.Range("C1").Value = "Email"
$ws2.Range("D1").Value = "Password"
$ws2.Range("E1").Value = "ConfirmPassword"

# Data row
$ws2.Range("C2").Value = "testingroro1@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "mailto:testingroro1@gmail.com")
$ws2.Range("D2").Value = "hello123"
$ws2.Range("E2").Value = "hello123"

$ws2.Range("A1").Value = "Firstname"
$ws2.Range("B1").Value = "Lastname "

$ws2.Range("A2").Value = "Munni"
$ws2.Range("B2").Value = "Roro"

$ws2.Rows.Item(2).RowHeight = 26

$ws2.Range("B2").Select()
